# Update the "GetsDebin" worksheet cells that contained the literal
# prefix "debin." before the field name "ori_trx_id" inside the JSON
# snippet text. The prefix is removed so the values read "ori_trx_id"
# instead of "debin.ori_trx_id".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GetsDebin")

$ws.Range("C5").Value = '"ori_trx_id":"9223000000000066168"'
$ws.Range("C6").Value = '"ori_trx_id":"12347"'
$ws.Range("C7").Value = '"ori_trx_id":"9223000000000"'
